# Update cryptos list - price (D) and volume (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched cells to stay text (matches the original inline-string
# cells) instead of being auto-coerced to numbers by values like "226.72".
$touchedRange = $ws.Range("D2:E51")
$touchedRange.NumberFormat = "@"

$updates = @(
    @{ Row = 2;  D = "34.159.72";  E = "  +0.49%  " },
    @{ Row = 3;  D = "1.791.02";   E = "  +0.45%  " },
    @{ Row = 4;  E = "  +0.15%  " },
    @{ Row = 5;  D = "226.72";     E = "  +0.67%  " },
    @{ Row = 6;  D = "0.548" },
    @{ Row = 7;  E = "  +0.13%  " },
    @{ Row = 8;  D = "31.92";      E = "  -1.22%  " },
    @{ Row = 9;  E = "  +1.14%  " },
    @{ Row = 10; D = "0.0691";     E = "  -2.02%  " },
    @{ Row = 11; D = "0.0947";     E = "  +1.12%  " },
    @{ Row = 12; D = "2.049.34";   E = "  +0.58%  " },
    @{ Row = 13; E = "  +1.34%  " },
    @{ Row = 14; D = "1.788.34";   E = "  +0.27%  " },
    @{ Row = 15; D = "34.110.94";  E = "  +0.36%  " },
    @{ Row = 16; E = "  -0.09%  " },
    @{ Row = 17; D = "4.19";       E = "  +0.75%  " },
    @{ Row = 18; E = "  +0.50%  " },
    @{ Row = 19; D = "245.61";     E = "  +0.77%  " },
    @{ Row = 20; E = "  -0.58%  " },
    @{ Row = 21; E = "  +0.05%  " },
    @{ Row = 22; E = "  +1.00%  " },
    @{ Row = 23; E = "  +0.26%  " },
    @{ Row = 24; E = "  -0.04%  " },
    @{ Row = 25; D = "161.10";     E = "  +0.60%  " },
    @{ Row = 26; E = "  +0.95%  " },
    @{ Row = 27; D = "16.36" },
    @{ Row = 28; E = "  +0.72%  " },
    @{ Row = 29; E = "  +0.25%  " },
    @{ Row = 30; E = "  -0.50%  " },
    @{ Row = 31; E = "  +1.38%  " },
    @{ Row = 32; E = "  +0.25%  " },
    @{ Row = 33; E = "  +2.87%  " },
    @{ Row = 34; E = "  -0.50%  " },
    @{ Row = 35; D = "1.460.69";   E = "  +4.59%  " },
    @{ Row = 36; D = "0.648";      E = "  -1.55%  " },
    @{ Row = 37; E = "  +7.52%  " },
    @{ Row = 38; E = "  +2.43%  " },
    @{ Row = 39; E = "  -0.15%  " },
    @{ Row = 40; D = "80.42";      E = "  +3.02%  " },
    @{ Row = 41; E = "  +0.48%  " },
    @{ Row = 42; E = "  +0.78%  " },
    @{ Row = 43; E = "  +0.28%  " },
    @{ Row = 44; D = "13.51";      E = "  +0.77%  " },
    @{ Row = 45; E = "  +2.91%  " },
    @{ Row = 46; D = "6.07";       E = "  +3.93%  " },
    @{ Row = 47; E = "  +0.08%  " },
    @{ Row = 48; E = "  -2.58%  " },
    @{ Row = 49; D = "1.950.41";   E = "  +0.95%  " },
    @{ Row = 50; D = "106.09";     E = "  -1.93%  " },
    @{ Row = 51; E = "  +0.11%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

# Restore the default ("Normal") style on the touched cells so no lingering
# per-cell style index is left behind (the source workbook has none here).
$touchedRange.Style = "Normal"
